$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "RawData"
$ws.Range("B2").Value = "Reports"
$ws.Range("B4").Value = "Parameters"
$ws.Range("A5").Value = "ReportsPrepared"
$ws.Range("B5").Value = "ReportsPrepared"
$ws.Range("A6").Value = "RoundFTE"
$ws.Range("B6").Value = 3

$ws.Range("B5").Select() | Out-Null
